$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ40713573",
    "summ40819175",
    "summ40924417",
    "summ41031199",
    "summ41143810",
    "summ41248147",
    "summ41350435",
    "summ41558600",
    "summ41676616",
    "summ41774684",
    "summ41875222",
    "summ41983025",
    "summ42085893",
    "summ42192982",
    "summ42297276",
    "summ42402819",
    "summ42522194",
    "summ42650720",
    "summ42782828",
    "summ42914716",
    "summ43056298",
    "summ43195600",
    "summ43333435",
    "summ43484067",
    "summ43625249",
    "summ43774252",
    "summ43906190",
    "summ44041923",
    "summ44196921",
    "summ44346256",
    "summ44482357",
    "summ44610407",
    "summ44747260",
    "summ44886319",
    "summ45023718",
    "summ45164389",
    "summ45317894",
    "summ45498045",
    "summ45634698",
    "summ45777170",
    "summ45905345",
    "summ46032383",
    "summ46173291",
    "summ46314275",
    "summ46440796",
    "summ46566413",
    "summ46712370",
    "summ46849030",
    "summ46975805",
    "summ47147524"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}

